$d = $word.ActiveDocument

# The paragraph currently holds its text as many <w:t> runs separated by
# manual line breaks (<w:br/>, represented in the Word text stream as the
# vertical-tab character). The edit collapses all of those line breaks into
# single spaces so the whole paragraph becomes one contiguous run of text.

# Step 1: turn every manual line break into a single space.
$d.Content.Find.Execute("^l", $false, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null

# Step 2: the very last line break (immediately before the paragraph mark)
# turned into a trailing space in step 1 - drop that stray space so the
# paragraph ends right after the final character, with no trailing break.
$d.Content.Find.Execute(" ^p", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
